$wb = $excel.ActiveWorkbook

# --- Add new "Assay" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

$data = @(
    @("Setting", "Value"),
    @("DMSO Tolerance", 0.005),
    @("Well Volume (µL)", 25),
    @("Backfill (µL)", 10),
    @("Allowed Error", 0.1),
    @("Destination Replicates", 1),
    @("Use Intermediate Plates", 1),
    @("DMSO Normalization", 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $assay.Cells.Item($row, 1).Value = $data[$i][0]
    $assay.Cells.Item($row, 2).Value = $data[$i][1]
}

[void]$assay.Range("A1:B8").Select()

# --- Compounds sheet: remove the (redundant) header style from row 1 ---
$compounds = $wb.Worksheets.Item("Compounds")
$compounds.Range("A1:F1").ClearFormats()

# --- Make "Patterns" the active sheet / selection, matching the saved view state ---
$patterns = $wb.Worksheets.Item("Patterns")
[void]$patterns.Activate()
[void]$patterns.Range("D12").Select()
